$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and 1h-volume-change (E) columns for rows 2-50 ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.357.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.376.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.377.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +3.92%  '
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("E12").Value = '  +2.40%  '
$ws.Range("E13").Value = '  +3.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("E15").Value = '  +4.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.808.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.178.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.379.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.76%  '
$ws.Range("E20").Value = '  +2.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.44%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.95%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.493.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '513.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0897'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.378'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '146.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.59%  '
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0526'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("E50").Value = '  +2.18%  '

# --- Rows 32 and 33 swap: Fetch.AI <-> Kaspa (with updated price/volume) ---
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.150'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.85%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.39'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.49%  '
